$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row2
$ws.Range("B2").Value = 0.1986531986531987
$ws.Range("C2").Value = 0.5353535353535354
$ws.Range("J2").Value = 0.0101010101010101
$ws.Range("P2").Value = 0.1481481481481481
$ws.Range("S2").Value = 0.1077441077441077

# row3
$ws.Range("B3").Value = 0.006060606060606061
$ws.Range("C3").Value = 0.0303030303030303
$ws.Range("J3").Value = 0.04242424242424243
$ws.Range("P3").Value = 0.7333333333333333
$ws.Range("S3").Value = 0.1878787878787879

# row6
$ws.Range("B6").Value = 0.04577464788732395
$ws.Range("D6").Value = 0.02112676056338028
$ws.Range("F6").Value = 0.09859154929577464
$ws.Range("J6").Value = 0.2640845070422535
$ws.Range("O6").Value = 0.01408450704225352
$ws.Range("Q6").Value = 0.09859154929577464
$ws.Range("R6").Value = 0.1056338028169014
$ws.Range("S6").Value = 0.352112676056338

# row7
$ws.Range("B7").Value = 0.08125
$ws.Range("D7").Value = 0.05
$ws.Range("F7").Value = 0.05625
$ws.Range("J7").Value = 0.1375
$ws.Range("O7").Value = 0.03125
$ws.Range("Q7").Value = 0.175
$ws.Range("R7").Value = 0.08749999999999999
$ws.Range("S7").Value = 0.38125

# row8
$ws.Range("B8").Value = 0.1056511056511057
$ws.Range("D8").Value = 0.02211302211302211
$ws.Range("F8").Value = 0.06388206388206388
$ws.Range("J8").Value = 0.1326781326781327
$ws.Range("O8").Value = 0.02457002457002457
$ws.Range("Q8").Value = 0.1646191646191646
$ws.Range("R8").Value = 0.085995085995086
$ws.Range("S8").Value = 0.4004914004914005

# row9
$ws.Range("B9").Value = 0.09016393442622951
$ws.Range("D9").Value = 0.01639344262295082
$ws.Range("E9").Value = 0.004098360655737705
$ws.Range("F9").Value = 0.06557377049180328
$ws.Range("J9").Value = 0.1311475409836066
$ws.Range("O9").Value = 0.01639344262295082
$ws.Range("Q9").Value = 0.1639344262295082
$ws.Range("R9").Value = 0.1065573770491803
$ws.Range("S9").Value = 0.4057377049180328

# row10
$ws.Range("B10").Value = 0.103988603988604
$ws.Range("D10").Value = 0.0170940170940171
$ws.Range("E10").Value = 0.001424501424501425
$ws.Range("F10").Value = 0.08475783475783476
$ws.Range("J10").Value = 0.1075498575498575
$ws.Range("O10").Value = 0.01353276353276353
$ws.Range("Q10").Value = 0.1994301994301994
$ws.Range("R10").Value = 0.09116809116809117
$ws.Range("S10").Value = 0.3810541310541311

# row11
$ws.Range("G11").Value = 0.1156716417910448
$ws.Range("J11").Value = 0.09701492537313433
$ws.Range("K11").Value = 0.1604477611940298
$ws.Range("L11").Value = 0.6156716417910447
$ws.Range("S11").Value = 0.01119402985074627

# row12
$ws.Range("G12").Value = 0.6987951807228916
$ws.Range("J12").Value = 0.2530120481927711
$ws.Range("K12").Value = 0.006024096385542169
$ws.Range("L12").Value = 0.01204819277108434
$ws.Range("S12").Value = 0.03012048192771084

# row13
$ws.Range("F13").Value = 0.06060606060606061
$ws.Range("G13").Value = 0.5454545454545454
$ws.Range("J13").Value = 0.3636363636363636
$ws.Range("S13").Value = 0.0303030303030303

# row15
$ws.Range("F15").Value = 0.0131578947368421
$ws.Range("H15").Value = 0.1743421052631579
$ws.Range("I15").Value = 0.09210526315789473
$ws.Range("J15").Value = 0.3881578947368421
$ws.Range("K15").Value = 0.06578947368421052
$ws.Range("M15").Value = 0.0131578947368421
$ws.Range("O15").Value = 0.1019736842105263
$ws.Range("S15").Value = 0.1513157894736842

# row16
$ws.Range("F16").Value = 0.01595744680851064
$ws.Range("H16").Value = 0.1968085106382979
$ws.Range("I16").Value = 0.0797872340425532
$ws.Range("J16").Value = 0.324468085106383
$ws.Range("K16").Value = 0.1223404255319149
$ws.Range("M16").Value = 0.01595744680851064
$ws.Range("O16").Value = 0.101063829787234
$ws.Range("S16").Value = 0.1436170212765958

# row17
$ws.Range("F17").Value = 0.02252252252252252
$ws.Range("H17").Value = 0.1373873873873874
$ws.Range("I17").Value = 0.1148648648648649
$ws.Range("J17").Value = 0.4279279279279279
$ws.Range("K17").Value = 0.09234234234234234
$ws.Range("M17").Value = 0.009009009009009009
$ws.Range("N17").Value = 0.002252252252252252
$ws.Range("O17").Value = 0.07657657657657657
$ws.Range("S17").Value = 0.1171171171171171

# row18
$ws.Range("F18").Value = 0.04741379310344827
$ws.Range("H18").Value = 0.1508620689655172
$ws.Range("I18").Value = 0.09913793103448276
$ws.Range("J18").Value = 0.4310344827586207
$ws.Range("K18").Value = 0.04741379310344827
$ws.Range("M18").Value = 0.01293103448275862
$ws.Range("O18").Value = 0.08189655172413793
$ws.Range("S18").Value = 0.1293103448275862

# row19
$ws.Range("F19").Value = 0.02386451116243264
$ws.Range("H19").Value = 0.1732101616628176
$ws.Range("I19").Value = 0.09930715935334873
$ws.Range("J19").Value = 0.4026173979984604
$ws.Range("K19").Value = 0.09468822170900693
$ws.Range("M19").Value = 0.01770592763664357
$ws.Range("O19").Value = 0.09699769053117784
$ws.Range("S19").Value = 0.09160892994611239

